$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: fill in remaining day-of-week "X" markers (H,J,L,N) ---
$ws.Cells.Item(16, 8).Value  = "X"   # H16
$ws.Cells.Item(16, 10).Value = "X"   # J16
$ws.Cells.Item(16, 12).Value = "X"   # L16
$ws.Cells.Item(16, 14).Value = "X"   # N16

# --- Row 17: fill in remaining day-of-week "X" markers (I,M) ---
$ws.Cells.Item(17, 9).Value  = "X"   # I17
$ws.Cells.Item(17, 13).Value = "X"   # M17

# --- Row 18 & 19: new trend_epi export rows ---
# Values are entered in the same order the original author typed them so
# that newly-created shared-string table entries line up with the source
# workbook (basic_CRF_extract.xlsx, trend_epi, keep_only_trend_epi_cols,
# basic_CRF_extract.sas7bdat, the CaseReportForm path, save_sas7bdat).
$ws.Cells.Item(18, 6).Value  = "basic_CRF_extract.xlsx"
$ws.Cells.Item(18, 1).Value  = "trend_epi"
$ws.Cells.Item(18, 3).Value  = "keep_only_trend_epi_cols"
$ws.Cells.Item(19, 6).Value  = "basic_CRF_extract.sas7bdat"
$ws.Cells.Item(18, 4).Value  = "//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/EPI SUMMARY/Trend analysis/_Current/_Source Data/CaseReportForm"
$ws.Cells.Item(19, 7).Value  = "save_sas7bdat"

$ws.Cells.Item(18, 2).Value  = "select * from all_cases;"
$ws.Cells.Item(18, 5).Value  = " "
$ws.Cells.Item(18, 8).Value  = "X"
$ws.Cells.Item(18, 9).Value  = "X"
$ws.Cells.Item(18, 10).Value = "X"
$ws.Cells.Item(18, 11).Value = "X"
$ws.Cells.Item(18, 12).Value = "X"
$ws.Cells.Item(18, 13).Value = "X"
$ws.Cells.Item(18, 14).Value = "X"

$ws.Cells.Item(19, 1).Value  = "trend_epi"
$ws.Cells.Item(19, 2).Value  = "select * from all_cases;"
$ws.Cells.Item(19, 3).Value  = "keep_only_trend_epi_cols"
$ws.Cells.Item(19, 4).Value  = "//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/EPI SUMMARY/Trend analysis/_Current/_Source Data/CaseReportForm"
$ws.Cells.Item(19, 5).Value  = " "
$ws.Cells.Item(19, 8).Value  = "X"
$ws.Cells.Item(19, 9).Value  = "X"
$ws.Cells.Item(19, 10).Value = "X"
$ws.Cells.Item(19, 11).Value = "X"
$ws.Cells.Item(19, 12).Value = "X"
$ws.Cells.Item(19, 13).Value = "X"
$ws.Cells.Item(19, 14).Value = "X"

# --- Column D needs to widen/autofit to accommodate the long new path text ---
$ws.Columns.Item(4).ColumnWidth = 216.6

# --- View state: scroll/select to match author's final view ---
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("H9").Select() | Out-Null
